$d = $word.ActiveDocument

# --- Change 1: "Hello, Soporte!" -> two runs "Hola" + ", Soporte!" -------
# (same bold/size formatting on both runs, just split into separate <w:r>
# elements as in the target compiled XML)
$p1 = $d.Paragraphs(1).Range
if ($p1.Text -eq "Hello, Soporte!" -or $p1.Text -eq "Hello, Soporte!`r") {
    $xml1 = '<w:wordDocument xmlns:w="http://schemas.microsoft.com/office/word/2003/wordml">' +
            '<w:body><w:p>' +
                '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Hola</w:t></w:r>' +
                '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">, Soporte!</w:t></w:r>' +
            '</w:p></w:body></w:wordDocument>'
    $p1.InsertXML($xml1)
}

# --- Change 2: drop w:hint="cs" from the "Al Tarikh" paragraph mark rPr --
$pTarikh = $null
foreach ($p in $d.Paragraphs) {
    $rpr = $p.Range.ParagraphFormat
    if ($p.Range.Text -eq "Prueba de Compilador DocxSerializer.`r" -or $p.Range.Text -eq "Prueba de Compilador DocxSerializer.") {
        $pTarikh = $p
    }
}
if ($pTarikh -ne $null) {
    $xml2 = '<w:wordDocument xmlns:w="http://schemas.microsoft.com/office/word/2003/wordml">' +
            '<w:body><w:p>' +
                '<w:pPr><w:rPr><w:rFonts w:ascii="Al Tarikh" w:hAnsi="Al Tarikh" w:cs="Al Tarikh"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Prueba de Compilador DocxSerializer.</w:t></w:r>' +
            '</w:p></w:body></w:wordDocument>'
    $pTarikh.Range.InsertXML($xml2)
}

Write-Output "done"
